$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 1: J1:N1 = FA, TA, SW, CW, QA (mirrors B1:F1)
$ws.Range("J1").Value = "FA"
$ws.Range("K1").Value = "TA"
$ws.Range("L1").Value = "SW"
$ws.Range("M1").Value = "CW"
$ws.Range("N1").Value = "QA"

# Category-average formulas for each contestant row (rows 2-12; 7/13 are spacer rows)
$ws.Range("J2").Formula = "=AVERAGE(B2,B14,B26)"
$ws.Range("K2").Formula = "=AVERAGE(C2,C14,C26)"
$ws.Range("L2").Formula = "=AVERAGE(D2,D14,D26)"
$ws.Range("M2").Formula = "=AVERAGE(E2,E14,E26)"
$ws.Range("N2").Formula = "=AVERAGE(F2,F14,F26)"

$ws.Range("J3").Formula = "=AVERAGE(B3,B15,B27)"
$ws.Range("K3").Formula = "=AVERAGE(C3,C15,C27)"
$ws.Range("L3").Formula = "=AVERAGE(D3,D15,D27)"
$ws.Range("M3").Formula = "=AVERAGE(E3,E15,E27)"
$ws.Range("N3").Formula = "=AVERAGE(F3,F15,F27)"

$ws.Range("J4").Formula = "=AVERAGE(B4,B16,B28)"
$ws.Range("K4").Formula = "=AVERAGE(C4,C16,C28)"
$ws.Range("L4").Formula = "=AVERAGE(D4,D16,D28)"
$ws.Range("M4").Formula = "=AVERAGE(E4,E16,E28)"
$ws.Range("N4").Formula = "=AVERAGE(F4,F16,F28)"

$ws.Range("J5").Formula = "=AVERAGE(B5,B17,B29)"
$ws.Range("K5").Formula = "=AVERAGE(C5,C17,C29)"
$ws.Range("L5").Formula = "=AVERAGE(D5,D17,D29)"
$ws.Range("M5").Formula = "=AVERAGE(E5,E17,E29)"
$ws.Range("N5").Formula = "=AVERAGE(F5,F17,F29)"

$ws.Range("J6").Formula = "=AVERAGE(B6,B18,B30)"
$ws.Range("K6").Formula = "=AVERAGE(C6,C18,C30)"
$ws.Range("L6").Formula = "=AVERAGE(D6,D18,D30)"
$ws.Range("M6").Formula = "=AVERAGE(E6,E18,E30)"
$ws.Range("N6").Formula = "=AVERAGE(F6,F18,F30)"

$ws.Range("J8").Formula = "=AVERAGE(B8,B20,B32)"
$ws.Range("K8").Formula = "=AVERAGE(C8,C20,C32)"
$ws.Range("L8").Formula = "=AVERAGE(D8,D20,D32)"
$ws.Range("M8").Formula = "=AVERAGE(E8,E20,E32)"
$ws.Range("N8").Formula = "=AVERAGE(F8,F20,F32)"

$ws.Range("J9").Formula = "=AVERAGE(B9,B21,B33)"
$ws.Range("K9").Formula = "=AVERAGE(C9,C21,C33)"
$ws.Range("L9").Formula = "=AVERAGE(D9,D21,D33)"
$ws.Range("M9").Formula = "=AVERAGE(E9,E21,E33)"
$ws.Range("N9").Formula = "=AVERAGE(F9,F21,F33)"

$ws.Range("J10").Formula = "=AVERAGE(B10,B22,B34)"
$ws.Range("K10").Formula = "=AVERAGE(C10,C22,C34)"
$ws.Range("L10").Formula = "=AVERAGE(D10,D22,D34)"
$ws.Range("M10").Formula = "=AVERAGE(E10,E22,E34)"
$ws.Range("N10").Formula = "=AVERAGE(F10,F22,F34)"

$ws.Range("J11").Formula = "=AVERAGE(B11,B23,B35)"
$ws.Range("K11").Formula = "=AVERAGE(C11,C23,C35)"
$ws.Range("L11").Formula = "=AVERAGE(D11,D23,D35)"
$ws.Range("M11").Formula = "=AVERAGE(E11,E23,E35)"
$ws.Range("N11").Formula = "=AVERAGE(F11,F23,F35)"

$ws.Range("J12").Formula = "=AVERAGE(B12,B24,B36)"
$ws.Range("K12").Formula = "=AVERAGE(C12,C24,C36)"
$ws.Range("L12").Formula = "=AVERAGE(D12,D24,D36)"
$ws.Range("M12").Formula = "=AVERAGE(E12,E24,E36)"
$ws.Range("N12").Formula = "=AVERAGE(F12,F24,F36)"

# Match the author's final selection/view state
$ws.Range("N14").Select()
